# Ran code for averaged intensities on spiral schemes.
# Re-ran the HKL-averaged-intensity pipeline with the Gaussian-Quadrature
# scheme plus three new spiral sampling schemes
# (Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space, Spiral-90deg-10rot-3space)
# inserted ahead of the existing NoRotation/Rotation/HexGrid rows, pushing those
# three rows down to the bottom of the table (rows 17-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of row 16 (A:M) down into the newly-needed rows 17-19
# so the new rows inherit the same bold/bordered "index" column formatting
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.998451683503045
$ws.Range("D10").Value = 4.564337571077123
$ws.Range("E10").Value = 0.5187705011526014
$ws.Range("F10").Value = 1.998451683503045
$ws.Range("G10").Value = 1.188758491135423
$ws.Range("H10").Value = 0.2971667452826968
$ws.Range("I10").Value = 0.6188239908525918
$ws.Range("J10").Value = 4.564337571077123
$ws.Range("K10").Value = 2.541554036114862
$ws.Range("L10").Value = 2.270002859808954
$ws.Range("M10").Value = 1.531051497167247

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9308434093319623
$ws.Range("D11").Value = 1.179769571137464
$ws.Range("E11").Value = 1.225287451067756
$ws.Range("F11").Value = 0.9308434093319623
$ws.Range("G11").Value = 0.5406257772082353
$ws.Range("H11").Value = 2.707851563680761
$ws.Range("I11").Value = 0.9134528147725448
$ws.Range("J11").Value = 1.179769571137464
$ws.Range("K11").Value = 1.20252851110261
$ws.Range("L11").Value = 1.066685960217286
$ws.Range("M11").Value = 1.249638431199787

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9147531115279629
$ws.Range("D12").Value = 1.173537554345956
$ws.Range("E12").Value = 1.228433984290763
$ws.Range("F12").Value = 0.9147531115279629
$ws.Range("G12").Value = 0.5429115775287642
$ws.Range("H12").Value = 2.712679033980858
$ws.Range("I12").Value = 0.9138772236164332
$ws.Range("J12").Value = 1.173537554345956
$ws.Range("K12").Value = 1.200985769318359
$ws.Range("L12").Value = 1.057869440423161
$ws.Range("M12").Value = 1.247698747548456

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9286515978401849
$ws.Range("D13").Value = 1.16919753509214
$ws.Range("E13").Value = 1.225918149930745
$ws.Range("F13").Value = 0.9286515978401849
$ws.Range("G13").Value = 0.5420923686349158
$ws.Range("H13").Value = 2.706548614336546
$ws.Range("I13").Value = 0.9132522747862823
$ws.Range("J13").Value = 1.16919753509214
$ws.Range("K13").Value = 1.197557842511442
$ws.Range("L13").Value = 1.063104720175813
$ws.Range("M13").Value = 1.247610090103469

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.007284000000000001
$ws.Range("D14").Value = 9.006064000000018
$ws.Range("E14").Value = 0.7151520000000015
$ws.Range("F14").Value = 0.007284000000000001
$ws.Range("G14").Value = 1.349999999999999
$ws.Range("H14").Value = 1.410384000000004
$ws.Range("I14").Value = 0.344304
$ws.Range("J14").Value = 9.006064000000018
$ws.Range("K14").Value = 4.86060800000001
$ws.Range("L14").Value = 2.433946000000005
$ws.Range("M14").Value = 2.138864666666671

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 15.29763750000003
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1.791962499999997
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 15.29763750000003
$ws.Range("K15").Value = 7.648818750000013
$ws.Range("L15").Value = 3.824409375000006
$ws.Range("M15").Value = 2.84826666666667

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.4308558880767964
$ws.Range("D16").Value = 9.066948037324806
$ws.Range("E16").Value = 0.3791769640960029
$ws.Range("F16").Value = 0.4308558880767964
$ws.Range("G16").Value = 1.5085157429248
$ws.Range("H16").Value = 0.4417423302655989
$ws.Range("I16").Value = 0.4029426081792019
$ws.Range("J16").Value = 9.066948037324806
$ws.Range("K16").Value = 4.723062500710404
$ws.Range("L16").Value = 2.576959194393601
$ws.Range("M16").Value = 2.038363595144534

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 1.002045427560552
$ws.Range("D17").Value = 1.018234851944184
$ws.Range("E17").Value = 1.00457824618631
$ws.Range("F17").Value = 1.002045427560552
$ws.Range("G17").Value = 0.9962642003881156
$ws.Range("H17").Value = 0.9966612880256277
$ws.Range("I17").Value = 1.005246049566393
$ws.Range("J17").Value = 1.018234851944184
$ws.Range("K17").Value = 1.011406549065247
$ws.Range("L17").Value = 1.006725988312899
$ws.Range("M17").Value = 1.003838343945197

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.051734273079805
$ws.Range("D18").Value = 0.7807726619362321
$ws.Range("E18").Value = 1.118007090715019
$ws.Range("F18").Value = 1.051734273079805
$ws.Range("G18").Value = 0.9652541115587777
$ws.Range("H18").Value = 1.046798381174745
$ws.Range("I18").Value = 0.9951049823057196
$ws.Range("J18").Value = 0.7807726619362321
$ws.Range("K18").Value = 0.9493898763256257
$ws.Range("L18").Value = 1.000562074702716
$ws.Range("M18").Value = 0.9929452501283832

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.251877343497905
$ws.Range("D19").Value = 0.7617457692782816
$ws.Range("E19").Value = 0.979378311364891
$ws.Range("F19").Value = 1.251877343497905
$ws.Range("G19").Value = 0.7166147418083997
$ws.Range("H19").Value = 1.415762595771091
$ws.Range("I19").Value = 1.081891772405074
$ws.Range("J19").Value = 0.7617457692782816
$ws.Range("K19").Value = 0.8705620403215864
$ws.Range("L19").Value = 1.061219691909746
$ws.Range("M19").Value = 1.03454508902094
